$d = $word.ActiveDocument

$replacements = @(
    @("23×12=", "74×61="),
    @("89×94=", "86×94="),
    @("48×23=", "35×92="),
    @("46×96=", "98×58="),
    @("82×89=", "65×54="),
    @("67×17=", "58×81="),
    @("98×27=", "53×96="),
    @("66×33=", "88×56="),
    @("21×30=", "39×42="),
    @("22×27=", "86×79="),
    @("80×81=", "13×95="),
    @("46×69=", "19×94="),
    @("45×84=", "82×40="),
    @("91×28=", "48×53="),
    @("36×11=", "30×25="),
    @("96×60=", "34×94="),
    @("37×78=", "84×78="),
    @("32×74=", "75×60="),
    @("24×67=", "33×68="),
    @("14×37=", "75×63="),
    @("98×82=", "84×69="),
    @("94×34=", "55×62="),
    @("60×81=", "53×23="),
    @("92×54=", "16×14="),
    @("36×18=", "25×60=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
